# Auto-generated PowerShell Excel COM-interop script applying the commit's changes.
#
# Summary:
#  - Sheets 1-4 ('Bidirectional A', 'D Lite', 'IDA', 'SMA'): translate the
#    'Metrica'/'Valor' header to English ('Metric'/'Value') and refresh the
#    benchmark numbers with the re-measured results.
#  - Sheet 5 ('RTAA (lookahead=25, move_limit=') is left untouched.
#  - Two new RTAA benchmark sheets are appended, mirroring sheet 5's layout:
#    'RTAA (L=25,M=3)' (still Portuguese headers) and
#    'RTAA (L=25, M=3)' (translated headers, larger bulk-modification count).

$wb = $excel.ActiveWorkbook

# Template sheet used to copy the header (bold+border) / data (centered) cell
# styles (style indexes 1 and 2 in styles.xml) onto the new sheets.
$styleTemplate = $wb.Worksheets.Item(1)

# --- Sheet 1: refresh header translation + benchmark values ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(1, 1).Value = "Metric"
$ws.Cells.Item(1, 2).Value = "Value"
$ws.Cells.Item(2, 2).Value = [double]"0.00028954200024600140750408172607421875"
$ws.Cells.Item(3, 2).Value = [double]"0.00016962500012596140532179778137589210018632002174854278564453125"
$ws.Cells.Item(4, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws.Cells.Item(5, 2).Value = [double]"4011.8000000000001818989403545856475830078125"
$ws.Cells.Item(6, 2).Value = [double]"0.00292205810546875"
$ws.Cells.Item(7, 2).Value = [double]"0.0028671264648437500867361737988403547205962240695953369140625"
$ws.Cells.Item(8, 2).Value = [double]"0.00652980804443359375"
$ws.Cells.Item(9, 2).Value = [double]"0.00652980804443359375"
$ws.Cells.Item(10, 2).Value = [double]"0.00004350000017439015209674835205078125"
$ws.Cells.Item(11, 2).Value = [double]"0.0000415419999626465074010588096609097874534199945628643035888671875"
$ws.Cells.Item(12, 2).Value = [double]"0.0000476091599557548769499314145203783255055896006524562835693359375"
$ws.Cells.Item(13, 2).Value = [double]"0.00002677832000699709061146414068588228474254719913005828857421875"
$ws.Cells.Item(14, 2).Value = 100

# --- Sheet 2: refresh header translation + benchmark values ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(1, 1).Value = "Metric"
$ws.Cells.Item(1, 2).Value = "Value"
$ws.Cells.Item(2, 2).Value = [double]"0.00014891700084262990476206833800887352481367997825145721435546875"
$ws.Cells.Item(3, 2).Value = [double]"0.000940791998800705187022686004638671875"
$ws.Cells.Item(4, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws.Cells.Item(5, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws.Cells.Item(6, 2).Value = [double]"0.0028533935546875"
$ws.Cells.Item(7, 2).Value = [double]"0.0028533935546875"
$ws.Cells.Item(8, 2).Value = [double]"0.07830810546875"
$ws.Cells.Item(9, 2).Value = [double]"0.0781616210937499944488848768742172978818416595458984375"
$ws.Cells.Item(10, 2).Value = [double]"0.00020233400027791501106337601623152977481367997825145721435546875"
$ws.Cells.Item(11, 2).Value = [double]"0.000250791999860666692256927490234375"
$ws.Cells.Item(12, 2).Value = [double]"0.0000461063100192404787974233271174995252295047976076602935791015625"
$ws.Cells.Item(13, 2).Value = [double]"0.000020923770061926920036353416687546769026084803044795989990234375"
$ws.Cells.Item(14, 2).Value = 100

# --- Sheet 3: refresh header translation + benchmark values ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(1, 1).Value = "Metric"
$ws.Cells.Item(1, 2).Value = "Value"
$ws.Cells.Item(2, 2).Value = [double]"0.0003603750010370276868343353271484375"
$ws.Cells.Item(3, 2).Value = [double]"0.061914374999105348251760005950927734375"
$ws.Cells.Item(4, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws.Cells.Item(5, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws.Cells.Item(6, 2).Value = [double]"0.0028533935546875"
$ws.Cells.Item(7, 2).Value = [double]"0.0028533935546875"
$ws.Cells.Item(8, 2).Value = [double]"0.00331878662109375"
$ws.Cells.Item(9, 2).Value = [double]"0.0031600952148437499132638262011596452794037759304046630859375"
$ws.Cells.Item(10, 2).Value = [double]"0.000042375000703032128512859344482421875"
$ws.Cells.Item(11, 2).Value = [double]"0.041044792000320740044116973876953125"
$ws.Cells.Item(12, 2).Value = [double]"0.0032698736900238149598318937449903387459926307201385498046875"
$ws.Cells.Item(13, 2).Value = [double]"0.0000328432900278130624816853855296727715540328063070774078369140625"
$ws.Cells.Item(14, 2).Value = 100

# --- Sheet 4: refresh header translation + benchmark values ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(1, 1).Value = "Metric"
$ws.Cells.Item(1, 2).Value = "Value"
$ws.Cells.Item(2, 2).Value = [double]"0.00037458299993886612355709075927734375"
$ws.Cells.Item(3, 2).Value = [double]"0.00006912499884492717683315277099609375"
$ws.Cells.Item(4, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws.Cells.Item(5, 2).Value = [double]"6223.6999999999998181010596454143524169921875"
$ws.Cells.Item(6, 2).Value = [double]"0.00292205810546875"
$ws.Cells.Item(7, 2).Value = [double]"0.002903747558593750173472347597680709441192448139190673828125"
$ws.Cells.Item(8, 2).Value = [double]"0.00208282470703125"
$ws.Cells.Item(9, 2).Value = [double]"0.00208282470703125"
$ws.Cells.Item(10, 2).Value = [double]"0.000044667000111076049506664276123046875"
$ws.Cells.Item(11, 2).Value = [double]"0.0000513750001118751174994207786550504124534199945628643035888671875"
$ws.Cells.Item(12, 2).Value = [double]"0.00004544505989542813148736322137466459025745280086994171142578125"
$ws.Cells.Item(13, 2).Value = [double]"0.0000302424900655751086758250834041206189795047976076602935791015625"
$ws.Cells.Item(14, 2).Value = 100

# --- Sheet 5 (RTAA lookahead=25 sheet) is left untouched ---

# --- New sheet: "RTAA (L=25,M=3)" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add($null, $lastSheet)
$ws6.Name = "RTAA (L=25,M=3)"
$ws6.Cells.Item(1, 1).Value = "Métrica"
$ws6.Cells.Item(1, 2).Value = "Valor"
$ws6.Cells.Item(2, 1).Value = "A* Time (s)"
$ws6.Cells.Item(3, 1).Value = "RTAA* Time (s)"
$ws6.Cells.Item(4, 1).Value = "A* Cost"
$ws6.Cells.Item(5, 1).Value = "RTAA* Cost"
$ws6.Cells.Item(6, 1).Value = "A* Peak Memory (MiB)"
$ws6.Cells.Item(7, 1).Value = "A* Avg Memory (MiB)"
$ws6.Cells.Item(8, 1).Value = "RTAA* Peak Memory (MiB)"
$ws6.Cells.Item(9, 1).Value = "RTAA* Avg Memory (MiB)"
$ws6.Cells.Item(10, 1).Value = "A* Recalc Time (s)"
$ws6.Cells.Item(11, 1).Value = "RTAA* Recalc Time (s)"
$ws6.Cells.Item(12, 1).Value = "RTAA* Bulk Avg Time (s)"
$ws6.Cells.Item(13, 1).Value = "A* Bulk Avg Time (s)"
$ws6.Cells.Item(14, 1).Value = "Bulk Modifications Count"
$ws6.Cells.Item(2, 2).Value = [double]"0.00019220799913455269215985243835831397518632002174854278564453125"
$ws6.Cells.Item(3, 2).Value = [double]"0.0002124579987139441072940826416015625"
$ws6.Cells.Item(4, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws6.Cells.Item(5, 2).Value = [double]"4781.100000000000363797880709171295166015625"
$ws6.Cells.Item(6, 2).Value = [double]"0.00290679931640625"
$ws6.Cells.Item(7, 2).Value = [double]"0.002864074707031249826527652402319290558807551860809326171875"
$ws6.Cells.Item(8, 2).Value = [double]"0.00658416748046875"
$ws6.Cells.Item(9, 2).Value = [double]"0.00658416748046875"
$ws6.Cells.Item(10, 2).Value = [double]"0.000044124999476480297744274139404296875"
$ws6.Cells.Item(11, 2).Value = [double]"0.00011341700155753639610127503356551414981367997825145721435546875"
$ws6.Cells.Item(12, 2).Value = [double]"0.000062579199948231689631938934326171875"
$ws6.Cells.Item(13, 2).Value = [double]"0.00002519559984648367051223404045057208122670999728143215179443359375"
$ws6.Cells.Item(14, 2).Value = 10
$styleTemplate.Range("A1:B1").Copy()
$ws6.Range("A1:B1").PasteSpecial(-4122)
$styleTemplate.Range("A2:B14").Copy()
$ws6.Range("A2:B14").PasteSpecial(-4122)
$ws6.Columns.Item(1).ColumnWidth = 25.16666666666667
$ws6.Columns.Item(2).ColumnWidth = 23.16666666666667
$ws6.PageSetup.LeftMargin = 0.75 * 72
$ws6.PageSetup.RightMargin = 0.75 * 72
$ws6.PageSetup.TopMargin = 1 * 72
$ws6.PageSetup.BottomMargin = 1 * 72
$ws6.PageSetup.HeaderMargin = 0.5 * 72
$ws6.PageSetup.FooterMargin = 0.5 * 72

# --- New sheet: "RTAA (L=25, M=3)" ---
$ws7 = $wb.Worksheets.Add($null, $ws6)
$ws7.Name = "RTAA (L=25, M=3)"
$ws7.Cells.Item(1, 1).Value = "Metric"
$ws7.Cells.Item(1, 2).Value = "Value"
$ws7.Cells.Item(2, 1).Value = "A* Time (s)"
$ws7.Cells.Item(3, 1).Value = "RTAA* Time (s)"
$ws7.Cells.Item(4, 1).Value = "A* Cost"
$ws7.Cells.Item(5, 1).Value = "RTAA* Cost"
$ws7.Cells.Item(6, 1).Value = "A* Peak Memory (MiB)"
$ws7.Cells.Item(7, 1).Value = "A* Avg Memory (MiB)"
$ws7.Cells.Item(8, 1).Value = "RTAA* Peak Memory (MiB)"
$ws7.Cells.Item(9, 1).Value = "RTAA* Avg Memory (MiB)"
$ws7.Cells.Item(10, 1).Value = "A* Recalc Time (s)"
$ws7.Cells.Item(11, 1).Value = "RTAA* Recalc Time (s)"
$ws7.Cells.Item(12, 1).Value = "RTAA* Bulk Avg Time (s)"
$ws7.Cells.Item(13, 1).Value = "A* Bulk Avg Time (s)"
$ws7.Cells.Item(14, 1).Value = "Bulk Modifications Count"
$ws7.Cells.Item(2, 2).Value = [double]"0.00015829200128791851075009400329207664981367997825145721435546875"
$ws7.Cells.Item(3, 2).Value = [double]"0.0002144169993698596954345703125"
$ws7.Cells.Item(4, 2).Value = [double]"3404.1999999999998181010596454143524169921875"
$ws7.Cells.Item(5, 2).Value = [double]"4781.100000000000363797880709171295166015625"
$ws7.Cells.Item(6, 2).Value = [double]"0.00290679931640625"
$ws7.Cells.Item(7, 2).Value = [double]"0.002864074707031249826527652402319290558807551860809326171875"
$ws7.Cells.Item(8, 2).Value = [double]"0.00658416748046875"
$ws7.Cells.Item(9, 2).Value = [double]"0.00658416748046875"
$ws7.Cells.Item(10, 2).Value = [double]"0.000044000000343658030033111572265625"
$ws7.Cells.Item(11, 2).Value = [double]"0.00011487499978102280674858992615128272518632002174854278564453125"
$ws7.Cells.Item(12, 2).Value = [double]"0.000089153259941667778042291570006483425459009595215320587158203125"
$ws7.Cells.Item(13, 2).Value = [double]"0.00003223035006158169965030657611038122922764159739017486572265625"
$ws7.Cells.Item(14, 2).Value = 100
$styleTemplate.Range("A1:B1").Copy()
$ws7.Range("A1:B1").PasteSpecial(-4122)
$styleTemplate.Range("A2:B14").Copy()
$ws7.Range("A2:B14").PasteSpecial(-4122)
$ws7.Columns.Item(1).ColumnWidth = 25.16666666666667
$ws7.Columns.Item(2).ColumnWidth = 23.16666666666667
$ws7.PageSetup.LeftMargin = 0.75 * 72
$ws7.PageSetup.RightMargin = 0.75 * 72
$ws7.PageSetup.TopMargin = 1 * 72
$ws7.PageSetup.BottomMargin = 1 * 72
$ws7.PageSetup.HeaderMargin = 0.5 * 72
$ws7.PageSetup.FooterMargin = 0.5 * 72

# Clear the clipboard marching-ants/copy mode left over from the style copies.
$excel.CutCopyMode = 0

